$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update hours (D12) and pause (F12) for row 12 - the H12 shared formula recalculates automatically
$ws.Range("D12").Value = 14
$ws.Range("F12").Value = 0.25

# Update the active cell selection to match the saved state
$ws.Range("O8").Select()
